$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.294.30"
$ws.Range("E2").Value = "'  +0.57%  "

# Row 3
$ws.Range("D3").Value = "'1.888.31"
$ws.Range("E3").Value = "'  -1.52%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'320.19"
$ws.Range("E5").Value = "'  -2.76%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.07%  "

# Row 7
$ws.Range("D7").Value = "'0.5057"
$ws.Range("E7").Value = "'  -2.99%  "

# Row 8
$ws.Range("D8").Value = "'0.4023"
$ws.Range("E8").Value = "'  -1.34%  "

# Row 9
$ws.Range("D9").Value = "'0.08296"
$ws.Range("E9").Value = "'  -2.61%  "

# Row 10
$ws.Range("E10").Value = "'  -1.66%  "

# Row 11
$ws.Range("D11").Value = "'1.104"
$ws.Range("E11").Value = "'  -1.66%  "

# Row 12
$ws.Range("D12").Value = "'24.07"
$ws.Range("E12").Value = "'  +7.83%  "

# Row 13
$ws.Range("D13").Value = "'1.889.13"
$ws.Range("E13").Value = "'  -1.41%  "

# Row 14
$ws.Range("D14").Value = "'6.351"
$ws.Range("E14").Value = "'  -1.25%  "

# Row 15
$ws.Range("D15").Value = "'7.271"
$ws.Range("E15").Value = "'  -1.56%  "

# Row 16
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "'  +0.10%  "

# Row 17
$ws.Range("D17").Value = "'92.70"
$ws.Range("E17").Value = "'  -2.48%  "

# Row 18
$ws.Range("D18").Value = "'0.00001101"
$ws.Range("E18").Value = "'  -1.10%  "

# Row 19
$ws.Range("D19").Value = "'0.06491"
$ws.Range("E19").Value = "'  -3.01%  "

# Row 20
$ws.Range("D20").Value = "'18.24"
$ws.Range("E20").Value = "'  -0.79%  "

# Row 21
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "'  -0.01%  "

# Row 22
$ws.Range("D22").Value = "'30.288.45"
$ws.Range("E22").Value = "'  +0.52%  "

# Row 23
$ws.Range("D23").Value = "'5.894"
$ws.Range("E23").Value = "'  -1.88%  "

# Row 24
$ws.Range("D24").Value = "'11.25"
$ws.Range("E24").Value = "'  -0.65%  "

# Row 25
$ws.Range("D25").Value = "'2.172"
$ws.Range("E25").Value = "'  -2.01%  "

# Row 26
$ws.Range("D26").Value = "'2.098.33"
$ws.Range("E26").Value = "'  -1.81%  "

# Row 27
$ws.Range("D27").Value = "'21.37"
$ws.Range("E27").Value = "'  +1.40%  "

# Row 28
$ws.Range("D28").Value = "'159.63"
$ws.Range("E28").Value = "'  -0.47%  "

# Row 29
$ws.Range("D29").Value = "'2.261"
$ws.Range("E29").Value = "'  -6.57%  "

# Row 30
$ws.Range("D30").Value = "'128.58"
$ws.Range("E30").Value = "'  -0.21%  "

# Row 31
$ws.Range("D31").Value = "'1.094"
$ws.Range("E31").Value = "'  +1.45%  "

# Row 32
$ws.Range("D32").Value = "'0.1039"
$ws.Range("E32").Value = "'  -2.16%  "

# Row 33
$ws.Range("D33").Value = "'5.990"
$ws.Range("E33").Value = "'  -0.71%  "

# Row 34
$ws.Range("D34").Value = "'3.722"
$ws.Range("E34").Value = "'  +2.42%  "

# Row 35
$ws.Range("D35").Value = "'0.02444"
$ws.Range("E35").Value = "'  -1.88%  "

# Row 36
$ws.Range("D36").Value = "'5.319"
$ws.Range("E36").Value = "'  +2.80%  "

# Row 37
$ws.Range("D37").Value = "'0.06427"
$ws.Range("E37").Value = "'  -2.71%  "

# Row 38
$ws.Range("D38").Value = "'0.2152"
$ws.Range("E38").Value = "'  -2.44%  "

# Row 39
$ws.Range("D39").Value = "'1.182"
$ws.Range("E39").Value = "'  -3.77%  "

# Row 40
$ws.Range("D40").Value = "'8.580"
$ws.Range("E40").Value = "'  -3.17%  "

# Row 41
$ws.Range("D41").Value = "'0.6367"
$ws.Range("E41").Value = "'  -2.57%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.38"
$ws.Range("E42").Value = "'  -2.20%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.217"
$ws.Range("E43").Value = "'  -2.21%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "'  +0.09%  "

# Row 45
$ws.Range("D45").Value = "'13.16"
$ws.Range("E45").Value = "'  -1.00%  "

# Row 46
$ws.Range("D46").Value = "'0.5961"
$ws.Range("E46").Value = "'  -2.96%  "

# Row 47
$ws.Range("D47").Value = "'2.118"
$ws.Range("E47").Value = "'  +2.05%  "

# Row 48
$ws.Range("D48").Value = "'3.639"
$ws.Range("E48").Value = "'  -2.84%  "

# Row 49
$ws.Range("D49").Value = "'123.55"
$ws.Range("E49").Value = "'  -0.63%  "

# Row 50
$ws.Range("D50").Value = "'1.213"
$ws.Range("E50").Value = "'  -2.43%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'78.54"
$ws.Range("E51").Value = "'  -1.43%  "
